$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Furps+
$ws.Range("A6").Value = "Furps+"
$ws.Range("C6").Value = 43963
$ws.Range("D6").Value = 0.375
$ws.Range("E6").Value = 0.4375

# Row 7: review af MockUp
$ws.Range("A7").Value = "review af MockUp"
$ws.Range("C7").Value = 43963
$ws.Range("D7").Value = 0.4375
$ws.Range("E7").Value = 0.45833333333333331

# Row 8: Risikoanalyse
$ws.Range("A8").Value = "Risikoanalyse"
$ws.Range("C8").Value = 43963
$ws.Range("D8").Value = 0.45833333333333331
$ws.Range("E8").Value = 0.58333333333333337

# Row 9: Review af Metrikker
$ws.Range("A9").Value = "Review af Metrikker"
$ws.Range("C9").Value = 43963
$ws.Range("D9").Value = 0.58333333333333337
$ws.Range("E9").Value = 0.60416666666666663

# Row 10: Fællesgennemgang af projektplan
$ws.Range("A10").Value = "Fællesgennemgang af projektplan"
$ws.Range("C10").Value = 43963
$ws.Range("D10").Value = 0.60416666666666663
$ws.Range("E10").Value = 0.6875

# Move the active cell selection to E12, matching the saved view state
$ws.Range("E12").Select() | Out-Null
